# MysensorsGW BOM sheet update
# Commit message: "Added RFM reset signal, and DIO5 for future use"
#
# Concrete changes (derived from the canonical OOXML diff):
#  - Refresh the KiCad export header (Date / Tool version strings).
#  - Mark several rows "NM" (Not Mounted) in the Vendor / Vendor part
#    columns (G/H) for parts that are now informational placeholders
#    (RFM reset + DIO5 test points, and the other unpopulated rows).
#  - Drop the now-unused Supplier1 ("Mouser" / "Itead recommended")
#    entries on those same rows (and a couple of others that lost their
#    supplier annotation).
#  - Move the sheet's active selection down to A29:A33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header refresh -------------------------------------------------
$ws.Range("B2").Value = "søn 03 jul 2016 22:33:14 CEST"
$ws.Range("B3").Value = "Eeschema 4.1.0-alpha+201606220817+6945~45~ubuntu16.04.1-product"

# --- Add "NM" markers to Vendor (G) / Vendor part (H) columns -------
$nmBothCols = @(7, 21, 22, 23, 25, 26, 27, 28)
foreach ($r in $nmBothCols) {
    $ws.Cells.Item($r, 7).Value = "NM"
    $ws.Cells.Item($r, 8).Value = "NM"
}

$nmColGOnly = @(35, 39, 40, 41, 44, 45, 46, 47, 48)
foreach ($r in $nmColGOnly) {
    $ws.Cells.Item($r, 7).Value = "NM"
}

# --- Remove stale Supplier1 (column I) entries -----------------------
$clearI = @(13, 22, 23, 26, 27, 28, 31, 34, 43)
foreach ($r in $clearI) {
    $ws.Cells.Item($r, 9).ClearContents()
}

# --- Update the active selection -------------------------------------
$ws.Range("A29:A33").Select()
